$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of data to append (dates stored as Excel serial numbers)
$newRows = @(
    @{ Row = 54; A = "Bruno Díaz"; B = 42882; C = 1; D = "Sprint 3 - Integración BackEnd y FrontEnd"; E = "Pruebas con servicios nuestros (Servicios del Cliente) para implementar soluciones genéricas" },
    @{ Row = 55; A = "Bruno Díaz"; B = 42883; C = 1; D = "Sprint 3 - FrontEnd"; E = "Eliminación de todo lo que tenía que ver con MARPH del proyecto del Front" },
    @{ Row = 56; A = "Bruno Díaz"; B = 42883; C = 2; D = "Sprint 3 - FrontEnd"; E = "Eliminación de todo aspecto que no fuera reponsivo (por si implementamos app mobile)" },
    @{ Row = 57; A = "Bruno Díaz"; B = 42883; C = 1; D = "Sprint 3 - Integración BackEnd y FrontEnd"; E = "Investigación sobre Token en .NET (No llegué a ninguna conclusión interesante… Vi alguna forma de crear nuestro propio token)" },
    @{ Row = 58; A = "Bruno Díaz"; B = 42883; C = 5; D = "Sprint 3 - Integración BackEnd y FrontEnd"; E = "Pruebas con servicios nuestros (Servicios del Cliente) para implementar soluciones genéricas" }
)

# Use an existing date-formatted cell as the format source so the new
# cells reuse the workbook's existing date style (numFmtId 14) instead
# of creating a brand-new number format.
$dateFormatSource = $ws.Cells.Item(53, 2)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $dateFormatSource.Copy() | Out-Null
    $ws.Cells.Item($r.Row, 2).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
$excel.CutCopyMode = $false

# Update view to reflect final selection/scroll position
$ws.Range("D57").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 32
